$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("README")
$ws2 = $wb.Worksheets.Item("Tabelle1")
$ws3 = $wb.Worksheets.Item("metabolites_in_mmol_L-1")

# --- Shared-string text updates (model renamed/rephrased) ---
# Old "BIOMASS_Ecoli_core_w_GAM" header (sheet2!G1 and sheet3!D1) becomes the
# new model id used throughout.
$ws2.Range("G1").Value = "BIOMASS_Ec_iML1515_core_75p37M"
$ws3.Range("D1").Value = "BIOMASS_Ec_iML1515_core_75p37M"

# README explanatory text: the old single line is replaced by two lines.
$ws1.Range("A1").Value = "metabolites concentration in mmol/L"
$ws1.Range("A2").Value = "biomass concentration in g/L"

# --- New fermentation data points (rows 13-27) on Tabelle1 ---
$ws2.Range("A13").Value = 5.5
$ws2.Range("B13").Value = 72.224953255494015
$ws2.Range("C13").Value = 5.6399241327975016
$ws2.Range("E13").Value = 15

$ws2.Range("A14").Value = 6
$ws2.Range("B14").Value = 66.275143181334016
$ws2.Range("C14").Value = 7.0607167421733008
$ws2.Range("D14").Value = 23.839999999999918
$ws2.Range("E14").Value = 15

$ws2.Range("A15").Value = 6.5
$ws2.Range("B15").Value = 58.480370071020012
$ws2.Range("C15").Value = 8.6985527400000002
$ws2.Range("D15").Value = 33.679999999999836
$ws2.Range("E15").Value = 15

$ws2.Range("A16").Value = 7
$ws2.Range("B16").Value = 47.992959960039002
$ws2.Range("C16").Value = 11.148210167463002
$ws2.Range("D16").Value = 43.029999999999973
$ws2.Range("E16").Value = 15

$ws2.Range("A17").Value = 7.5
$ws2.Range("B17").Value = 33.998641326399607
$ws2.Range("C17").Value = 15.568524718173004
$ws2.Range("D17").Value = 53.009999999999991
$ws2.Range("E17").Value = 15

$ws2.Range("A18").Value = 8
$ws2.Range("B18").Value = 9.2437780257432021
$ws2.Range("C18").Value = 11.493629696768402
$ws2.Range("D18").Value = 76.309999999999945
$ws2.Range("E18").Value = 15

$ws2.Range("A19").Value = 8.5
$ws2.Range("B19").Value = 0.18628907533511402
$ws2.Range("C19").Value = 10.302101942443201
$ws2.Range("D19").Value = 90.330000000000155
$ws2.Range("E19").Value = 15

$ws2.Range("A20").Value = 9
$ws2.Range("B20").Value = 0.20577774782144703
$ws2.Range("C20").Value = 2.1002365638972007
$ws2.Range("E20").Value = 15

$ws2.Range("A21").Value = 9.5
$ws2.Range("B21").Value = 0.17025155364837602
$ws2.Range("C21").Value = 0.05675834658025801
$ws2.Range("E21").Value = 15

$ws2.Range("A22").Value = 10
$ws2.Range("B22").Value = 0.19101586889403005
$ws2.Range("C22").Value = 0.03764762621047801
$ws2.Range("D22").Value = 88.479999999999791
$ws2.Range("E22").Value = 15

$ws2.Range("A23").Value = 10.5
$ws2.Range("B23").Value = 0.19088147625419705
$ws2.Range("C23").Value = 0.03420560889126001
$ws2.Range("E23").Value = 15

$ws2.Range("A24").Value = 11
$ws2.Range("B24").Value = 0.22064748880284002
$ws2.Range("C24").Value = 0.038400775901883004
$ws2.Range("E24").Value = 15

$ws2.Range("A25").Value = 11.5
$ws2.Range("B25").Value = 0.24528787914943803
$ws2.Range("C25").Value = 0.03407034639615301
$ws2.Range("D25").Value = 87.410000000000082
$ws2.Range("E25").Value = 15

$ws2.Range("A26").Value = 12
$ws2.Range("B26").Value = 0.26730435106201506
$ws2.Range("C26").Value = 0.036958120929954
$ws2.Range("E26").Value = 5

$ws2.Range("A27").Value = 12.5
$ws2.Range("B27").Value = 0.27065111922873003
$ws2.Range("C27").Value = 0.03812749636996801
$ws2.Range("D27").Value = 82.369999999999891
$ws2.Range("E27").Value = 15

# Extend the c_biomass formula (D/E) down through the new rows - only for the
# rows that actually have a biomass sample (D filled in), mirroring the
# original fill-down pattern. Row 26 is the one exception: its D is blank but
# the formula was still dragged into F26 (evaluating to 0).
$ws2.Range("F14").Formula = "=D14/E14"
$ws2.Range("F15").Formula = "=D15/E15"
$ws2.Range("F16").Formula = "=D16/E16"
$ws2.Range("F17").Formula = "=D17/E17"
$ws2.Range("F18").Formula = "=D18/E18"
$ws2.Range("F19").Formula = "=D19/E19"
$ws2.Range("F22").Formula = "=D22/E22"
$ws2.Range("F25").Formula = "=D25/E25"
$ws2.Range("F26").Formula = "=D26/E26"
$ws2.Range("F27").Formula = "=D27/E27"

# G mirrors F where a c_biomass value was actually computed (row 26's formula
# evaluates to 0 but, like the source file, gets no G mirror value).
$ws2.Range("G14").Value = 1.5893333333333279
$ws2.Range("G15").Value = 2.2453333333333223
$ws2.Range("G16").Value = 2.8686666666666647
$ws2.Range("G17").Value = 3.5339999999999994
$ws2.Range("G18").Value = 5.0873333333333299
$ws2.Range("G19").Value = 6.02200000000001
$ws2.Range("G22").Value = 5.898666666666653
$ws2.Range("G25").Value = 5.827333333333339
$ws2.Range("G27").Value = 5.4913333333333263

# --- Mirror the same new rows on the metabolites_in_mmol_L-1 sheet (A-C same,
#     D mirrors Tabelle1's G column) ---
$ws3.Range("A13").Value = 5.5
$ws3.Range("B13").Value = 72.224953255494015
$ws3.Range("C13").Value = 5.6399241327975016

$ws3.Range("A14").Value = 6
$ws3.Range("B14").Value = 66.275143181334016
$ws3.Range("C14").Value = 7.0607167421733008
$ws3.Range("D14").Value = 1.5893333333333279

$ws3.Range("A15").Value = 6.5
$ws3.Range("B15").Value = 58.480370071020012
$ws3.Range("C15").Value = 8.6985527400000002
$ws3.Range("D15").Value = 2.2453333333333223

$ws3.Range("A16").Value = 7
$ws3.Range("B16").Value = 47.992959960039002
$ws3.Range("C16").Value = 11.148210167463002
$ws3.Range("D16").Value = 2.8686666666666647

$ws3.Range("A17").Value = 7.5
$ws3.Range("B17").Value = 33.998641326399607
$ws3.Range("C17").Value = 15.568524718173004
$ws3.Range("D17").Value = 3.5339999999999994

$ws3.Range("A18").Value = 8
$ws3.Range("B18").Value = 9.2437780257432021
$ws3.Range("C18").Value = 11.493629696768402
$ws3.Range("D18").Value = 5.0873333333333299

$ws3.Range("A19").Value = 8.5
$ws3.Range("B19").Value = 0.18628907533511402
$ws3.Range("C19").Value = 10.302101942443201
$ws3.Range("D19").Value = 6.02200000000001

$ws3.Range("A20").Value = 9
$ws3.Range("B20").Value = 0.20577774782144703
$ws3.Range("C20").Value = 2.1002365638972007

$ws3.Range("A21").Value = 9.5
$ws3.Range("B21").Value = 0.17025155364837602
$ws3.Range("C21").Value = 0.05675834658025801

$ws3.Range("A22").Value = 10
$ws3.Range("B22").Value = 0.19101586889403005
$ws3.Range("C22").Value = 0.03764762621047801
$ws3.Range("D22").Value = 5.898666666666653

$ws3.Range("A23").Value = 10.5
$ws3.Range("B23").Value = 0.19088147625419705
$ws3.Range("C23").Value = 0.03420560889126001

$ws3.Range("A24").Value = 11
$ws3.Range("B24").Value = 0.22064748880284002
$ws3.Range("C24").Value = 0.038400775901883004

$ws3.Range("A25").Value = 11.5
$ws3.Range("B25").Value = 0.24528787914943803
$ws3.Range("C25").Value = 0.03407034639615301
$ws3.Range("D25").Value = 5.827333333333339

$ws3.Range("A26").Value = 12
$ws3.Range("B26").Value = 0.26730435106201506
$ws3.Range("C26").Value = 0.036958120929954

$ws3.Range("A27").Value = 12.5
$ws3.Range("B27").Value = 0.27065111922873003
$ws3.Range("C27").Value = 0.03812749636996801
$ws3.Range("D27").Value = 5.4913333333333263

# --- Selection / active-sheet bookkeeping to mirror the final saved state ---
# README: cursor left on A2 (last edited cell).
$ws1.Range("A2").Select()
# metabolites sheet: cursor left on D26.
$ws3.Range("D26").Select()
# Tabelle1 ends up the active tab, cursor on G26.
$ws2.Range("G26").Select()
